$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 298, shifting existing rows 298:392 down to 299:393
$ws.Rows.Item(298).Insert()

# Populate the newly inserted row 298 with the new record's data
$ws.Cells.Item(298, 1).Value = 4
$ws.Cells.Item(298, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(298, 3).Value = "Los Lagos"
$ws.Cells.Item(298, 4).Value = 45093
$ws.Cells.Item(298, 5).Value = 10
$ws.Cells.Item(298, 6).Value = "Fruta"
$ws.Cells.Item(298, 7).Value = 100108
$ws.Cells.Item(298, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(298, 9).Value = 100108002
$ws.Cells.Item(298, 10).Value = "Mango"
$ws.Cells.Item(298, 11).Value = "Sin especificar"
$ws.Cells.Item(298, 12).Value = "Primera"
$ws.Cells.Item(298, 13).Value = 200
$ws.Cells.Item(298, 14).Value = 8500
$ws.Cells.Item(298, 15).Value = 9000
$ws.Cells.Item(298, 16).Value = 8750
$ws.Cells.Item(298, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(298, 18).Value = "Perú"
$ws.Cells.Item(298, 19).Value = 2188
$ws.Cells.Item(298, 20).Value = 4

# Preserve the date number format on column D for the new row (style used by the rest of column D)
$ws.Cells.Item(298, 4).NumberFormat = $ws.Cells.Item(299, 4).NumberFormat
